$d = $word.ActiveDocument

$replacements = @(
    @{old="41×32=1312"; new="97×86=8342"},
    @{old="32×81=2592"; new="29×94=2726"},
    @{old="36×33=1188"; new="31×39=1209"},
    @{old="91×71=6461"; new="12×63=756"},
    @{old="44×95=4180"; new="19×17=323"},
    @{old="66×63=4158"; new="73×48=3504"},
    @{old="56×26=1456"; new="16×87=1392"},
    @{old="91×97=8827"; new="91×38=3458"},
    @{old="26×34=884"; new="16×88=1408"},
    @{old="53×37=1961"; new="32×43=1376"},
    @{old="61×13=793"; new="47×45=2115"},
    @{old="79×44=3476"; new="14×75=1050"},
    @{old="16×85=1360"; new="56×97=5432"},
    @{old="99×95=9405"; new="38×24=912"},
    @{old="27×81=2187"; new="77×84=6468"},
    @{old="72×53=3816"; new="20×35=700"},
    @{old="13×89=1157"; new="84×69=5796"},
    @{old="83×89=7387"; new="54×95=5130"},
    @{old="84×91=7644"; new="23×47=1081"},
    @{old="78×34=2652"; new="57×85=4845"},
    @{old="11×49=539"; new="62×21=1302"},
    @{old="81×31=2511"; new="34×52=1768"},
    @{old="53×60=3180"; new="63×13=819"},
    @{old="17×85=1445"; new="73×31=2263"},
    @{old="30×65=1950"; new="62×34=2108"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
